$wb = $excel.ActiveWorkbook

# A new handoff xliff was generated for file "003e9f1a-59fd-4961-9cec-bd93a64528b4"
# (row 5 of every sheet). Update the recorded timestamps to reflect the new
# handoff report generation.

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G5").Value = "2016-09-07 17:03:50"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H5").Value = "2016-09-07 17:03:44"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H5").Value = "2016-09-07 17:03:50"
